$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Level LC Description French" cell (F4) is being removed — its
# translation will now live in AppTranslations instead of the xlsx fixture.
$ws.Range("F4").ClearContents()
